$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> FAPs (target cluster), refreshed TPM-derived metrics
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.28161
$ws.Range("H2").Value = 0.84483
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.6153926666666667
$ws.Range("N2").Value = 1.846178
$ws.Range("O2").Value = 0.172697186719763
$ws.Range("P2").Value = 0.172697186719763
$ws.Range("Q2").Value = 0.17330072886
$ws.Range("R2").Value = 1.55970655974
$ws.Range("S2").Value = 0.172697186719763
$ws.Range("T2").Value = 0.172697186719763

# Row 3: FAPs -> MuSCs (target cluster), refreshed TPM-derived metrics
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.28161
$ws.Range("H3").Value = 0.84483
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.433374
$ws.Range("N3").Value = 1.300122
$ws.Range("O3").Value = 0.1216174235596306
$ws.Range("P3").Value = 0.1216174235596306
$ws.Range("Q3").Value = 0.12204245214
$ws.Range("R3").Value = 1.09838206926
$ws.Range("S3").Value = 0.1216174235596306
$ws.Range("T3").Value = 0.1216174235596306

# Row 4: MuSCs -> Resolving-Mac (target cluster), refreshed TPM-derived metrics
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.28161
$ws.Range("H4").Value = 0.84483
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.514653666666667
$ws.Range("N4").Value = 7.543961
$ws.Range("O4").Value = 0.7056853897206065
$ws.Range("P4").Value = 0.7056853897206065
$ws.Range("Q4").Value = 0.70815161907
$ws.Range("R4").Value = 6.37336457163
$ws.Range("S4").Value = 0.7056853897206065
$ws.Range("T4").Value = 0.7056853897206065

# Row 5 (previously Resolving-Mac target) no longer exists in the refreshed output
$ws.Rows("5:5").Delete()
